# Add two new columns, I (I0) and J (IF), to the weekly data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1), styled to match the other header cells (bold, thin
# box border, centered/top aligned) ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("I1:J1").Font.Bold = $true
$ws.Range("I1:J1").HorizontalAlignment = -4108
$ws.Range("I1:J1").VerticalAlignment = -4160
$ws.Range("I1:J1").Borders.LineStyle = 1

# --- Data rows (2-55) ---
$rowNums = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51,52,53,54,55)
$iVals   = @(8,9,8,8,9,9,9,9,9,9,9,9,10,9,9,10,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,10,9,9,9,9,9,9,9,9,9,8,6,6,4)
$jVals   = @(8,9,8,9,9,9,9,9,9,9,9,9,10,9,9,10,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,10,9,9,9,9,9,9,9,9,9,8,6,6,4)

for ($k = 0; $k -lt $rowNums.Length; $k++) {
    $r = $rowNums[$k]
    $ws.Cells.Item($r, 9).Value = $iVals[$k]
    $ws.Cells.Item($r, 10).Value = $jVals[$k]
}

Write-Output "Added I0/IF columns across rows 1-55"
